# Insert a new price record at row 195 on the data sheet ("Hortaliza, Vega
# Monumental Concepción - Zapallo italiano" subset). Inserting the row shifts
# the previously existing rows 195-208 down to 196-209 (their contents are
# left untouched by the Insert call), and the newly opened row 195 is filled
# in with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 195..208 down to 196..209, opening up a blank row 195.
$ws.Rows.Item(195).Insert()

# Populate the new row 195 with the new record's data.
$ws.Range("A195").Value = 11
$ws.Range("B195").Value = "Vega Monumental Concepción"
$ws.Range("C195").Value = "Bíobío"
$ws.Range("D195").Value = 45013
$ws.Range("E195").Value = 8
$ws.Range("F195").Value = 100112032
$ws.Range("G195").Value = "Zapallo italiano"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 220
$ws.Range("K195").Value = 6000
$ws.Range("L195").Value = 6500
$ws.Range("M195").Value = 6227
$ws.Range("N195").Value = "$/caja 50 unidades"
$ws.Range("O195").Value = "Región Metropolitana"
$ws.Range("P195").Value = 125
$ws.Range("Q195").Value = 50
$ws.Range("R195").Value = "Hortaliza"
